$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '30.308.13'
$ws.Range("D3").Value = '1.929.41'
$ws.Range("E4").Value = '  +0.27%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.7426'
$ws.Range("E5").Value = '  +3.12%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '243.70'
$ws.Range("E6").Value = '  -2.32%  '
$ws.Range("E7").Value = '  +0.30%  '
$ws.Range("B8").Value = 'Solana'
$ws.Range("C8").Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '27.48'
$ws.Range("E8").Value = '  -1.84%  '
$ws.Range("B9").Value = 'Cardano'
$ws.Range("C9").Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.3149'
$ws.Range("E9").Value = '  -1.72%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.07040'
$ws.Range("E10").Value = '  -0.84%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.08039'
$ws.Range("E11").Value = '  +0.38%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.7765'
$ws.Range("E12").Value = '  -1.52%  '
$ws.Range("D13").Value = '1.938.36'
$ws.Range("E13").Value = '  +0.47%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.368'
$ws.Range("E14").Value = '  -0.18%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '92.98'
$ws.Range("E15").Value = '  -1.83%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '14.47'
$ws.Range("E16").Value = '  -1.33%  '
$ws.Range("D17").Value = '30.318.91'
$ws.Range("E17").Value = '  +0.06%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '5.963'
$ws.Range("E18").Value = '  +4.02%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '250.01'
$ws.Range("E19").Value = '  -2.72%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.000007948'
$ws.Range("E20").Value = '  -1.52%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '1.003'
$ws.Range("E21").Value = '  +0.31%  '
$ws.Range("D22").Value = '2.160.67'
$ws.Range("E22").Value = '  -0.96%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '1.003'
$ws.Range("E23").Value = '  +0.28%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '6.648'
$ws.Range("E24").Value = '  -2.47%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '9.556'
$ws.Range("E25").Value = '  -0.10%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '165.69'
$ws.Range("E26").Value = '  +0.67%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '19.01'
$ws.Range("E27").Value = '  -0.41%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.1287'
$ws.Range("E28").Value = '  +0.38%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.164'
$ws.Range("E29").Value = '  -5.56%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.563'
$ws.Range("E30").Value = '  +2.02%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.356'
$ws.Range("E31").Value = '  -0.27%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.413'
$ws.Range("E32").Value = '  -0.10%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.110'
$ws.Range("E33").Value = '  -0.91%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.05214'
$ws.Range("E34").Value = '  +1.89%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.309'
$ws.Range("E35").Value = '  +1.69%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.7544'
$ws.Range("E36").Value = '  +0.76%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.768'
$ws.Range("E37").Value = '  +0.00%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.01951'
$ws.Range("E38").Value = '  -1.87%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.791'
$ws.Range("E39").Value = '  -0.20%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '6.520'
$ws.Range("E40").Value = '  +1.89%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '76.69'
$ws.Range("E41").Value = '  -2.02%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.4499'
$ws.Range("E42").Value = '  -0.56%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.951'
$ws.Range("E43").Value = '  -2.11%  '
$ws.Range("B44").Value = 'PaxDollar'
$ws.Range("C44").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.003'
$ws.Range("E44").Value = '  +0.30%  '
$ws.Range("B45").Value = 'TrustWalletToken'
$ws.Range("C45").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.8418'
$ws.Range("E45").Value = '  -0.49%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '7.675'
$ws.Range("E46").Value = '  +2.51%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '9.948'
$ws.Range("E47").Value = '  +1.15%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '101.29'
$ws.Range("E48").Value = '  +0.22%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '37.56'
$ws.Range("E49").Value = '  +1.92%  '
$ws.Range("E50").Value = '  -1.65%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.1233'
$ws.Range("E51").Value = '  +7.78%  '
